# Edit script: insert two new grape price rows (rows 345 and 346)
# into the "Fruta, Feria Lagunitas de Puerto Montt - Uva" sheet,
# pushing the existing data (old rows 345..394) down to rows 347..396.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 345. This shifts every
# row that was at 345 or below down by two, which is exactly the
# behaviour required by the diff (old row 345 -> new row 347, ...,
# old row 394 -> new row 396).
$ws.Rows("345:346").Insert()

# --- Populate new row 345 ---
$ws.Cells.Item(345, 1).Value = 4
$ws.Cells.Item(345, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(345, 3).Value = "Los Lagos"
$ws.Cells.Item(345, 4).Value = 45258
$ws.Cells.Item(345, 5).Value = 10
$ws.Cells.Item(345, 6).Value = "Fruta"
$ws.Cells.Item(345, 7).Value = 100109
$ws.Cells.Item(345, 8).Value = "Uva"
$ws.Cells.Item(345, 9).Value = 100109001
$ws.Cells.Item(345, 10).Value = "Uva"
$ws.Cells.Item(345, 11).Value = "Red Globe"
$ws.Cells.Item(345, 12).Value = "Primera"
$ws.Cells.Item(345, 13).Value = 200
$ws.Cells.Item(345, 14).Value = 28000
$ws.Cells.Item(345, 15).Value = 28000
$ws.Cells.Item(345, 16).Value = 28000
$ws.Cells.Item(345, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(345, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(345, 19).Value = 2800
$ws.Cells.Item(345, 20).Value = 10

# --- Populate new row 346 ---
$ws.Cells.Item(346, 1).Value = 4
$ws.Cells.Item(346, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(346, 3).Value = "Los Lagos"
$ws.Cells.Item(346, 4).Value = 45258
$ws.Cells.Item(346, 5).Value = 10
$ws.Cells.Item(346, 6).Value = "Fruta"
$ws.Cells.Item(346, 7).Value = 100109
$ws.Cells.Item(346, 8).Value = "Uva"
$ws.Cells.Item(346, 9).Value = 100109001
$ws.Cells.Item(346, 10).Value = "Uva"
$ws.Cells.Item(346, 11).Value = "Superior Seedless"
$ws.Cells.Item(346, 12).Value = "Primera"
$ws.Cells.Item(346, 13).Value = 200
$ws.Cells.Item(346, 14).Value = 28000
$ws.Cells.Item(346, 15).Value = 28000
$ws.Cells.Item(346, 16).Value = 28000
$ws.Cells.Item(346, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(346, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(346, 19).Value = 2800
$ws.Cells.Item(346, 20).Value = 10
